$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (default/no explicit formatting) used to keep written
# cells as plain text without leaving a stray number-format style behind.
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "'62.605.10"
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Value = "'  -1.24%  "
$ws.Range("E2").Style = $defaultStyle

$ws.Range("D3").Value = "'3.015.28"
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Value = "'  -1.59%  "
$ws.Range("E3").Style = $defaultStyle

$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("E4").Style = $defaultStyle

$ws.Range("D5").Value = "'584.71"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "'  -0.57%  "
$ws.Range("E5").Style = $defaultStyle

$ws.Range("D6").Value = "'147.57"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "'  -4.65%  "
$ws.Range("E6").Style = $defaultStyle

$ws.Range("E7").Value = "'  -0.06%  "
$ws.Range("E7").Style = $defaultStyle

$ws.Range("E8").Value = "'  -3.08%  "
$ws.Range("E8").Style = $defaultStyle

$ws.Range("D9").Value = "'3.015.66"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "'  -1.50%  "
$ws.Range("E9").Style = $defaultStyle

$ws.Range("E10").Value = "'  -3.41%  "
$ws.Range("E10").Style = $defaultStyle

$ws.Range("D11").Value = "'5.79"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "'  -0.57%  "
$ws.Range("E11").Style = $defaultStyle

$ws.Range("E12").Value = "'  -1.82%  "
$ws.Range("E12").Style = $defaultStyle

$ws.Range("E13").Value = "'  -2.72%  "
$ws.Range("E13").Style = $defaultStyle

$ws.Range("D14").Value = "'35.01"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "'  -5.22%  "
$ws.Range("E14").Style = $defaultStyle

$ws.Range("E15").Value = "'  +2.48%  "
$ws.Range("E15").Style = $defaultStyle

$ws.Range("D16").Value = "'3.516.82"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "'  -1.49%  "
$ws.Range("E16").Style = $defaultStyle

$ws.Range("D17").Value = "'7.04"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "'  -1.14%  "
$ws.Range("E17").Style = $defaultStyle

$ws.Range("D18").Value = "'62.583.20"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = "'  -1.33%  "
$ws.Range("E18").Style = $defaultStyle

$ws.Range("D19").Value = "'3.015.92"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "'  -1.53%  "
$ws.Range("E19").Style = $defaultStyle

$ws.Range("D20").Value = "'465.97"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "'  -1.21%  "
$ws.Range("E20").Style = $defaultStyle

$ws.Range("E21").Value = "'  -2.55%  "
$ws.Range("E21").Style = $defaultStyle

$ws.Range("E22").Value = "'  -2.44%  "
$ws.Range("E22").Style = $defaultStyle

$ws.Range("D24").Value = "'2.34"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "'  -3.36%  "
$ws.Range("E24").Style = $defaultStyle

$ws.Range("D25").Value = "'80.30"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "'  -0.39%  "
$ws.Range("E25").Style = $defaultStyle

$ws.Range("D26").Value = "'12.43"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "'  -2.81%  "
$ws.Range("E26").Style = $defaultStyle

$ws.Range("D27").Value = "'10.29"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "'  -0.74%  "
$ws.Range("E27").Style = $defaultStyle

$ws.Range("E28").Value = "'  +0.13%  "
$ws.Range("E28").Style = $defaultStyle

$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "'  -0.12%  "
$ws.Range("E29").Style = $defaultStyle

$ws.Range("E30").Value = "'  -1.01%  "
$ws.Range("E30").Style = $defaultStyle

$ws.Range("D31").Value = "'7.16"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "'  -4.28%  "
$ws.Range("E31").Style = $defaultStyle

$ws.Range("E32").Value = "'  -0.38%  "
$ws.Range("E32").Style = $defaultStyle

$ws.Range("D33").Value = "'27.62"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "'  +1.83%  "
$ws.Range("E33").Style = $defaultStyle

$ws.Range("E34").Value = "'  -4.22%  "
$ws.Range("E34").Style = $defaultStyle

$ws.Range("D35").Value = "'1.03"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "'  -0.46%  "
$ws.Range("E35").Style = $defaultStyle

$ws.Range("D36").Value = "'0.0₃0797"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "'  -2.75%  "
$ws.Range("E36").Style = $defaultStyle

$ws.Range("E37").Value = "'  -3.85%  "
$ws.Range("E37").Style = $defaultStyle

$ws.Range("D38").Value = "'2.14"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "'  -3.01%  "
$ws.Range("E38").Style = $defaultStyle

$ws.Range("D39").Value = "'50.37"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "'  -0.56%  "
$ws.Range("E39").Style = $defaultStyle

$ws.Range("D40").Value = "'8.98"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "'  -2.41%  "
$ws.Range("E40").Style = $defaultStyle

$ws.Range("D41").Value = "'2.94"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "'  -11.25%  "
$ws.Range("E41").Style = $defaultStyle

$ws.Range("D42").Value = "'423.18"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "'  -3.17%  "
$ws.Range("E42").Style = $defaultStyle

$ws.Range("D43").Value = "'0.113"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "'  +1.65%  "
$ws.Range("E43").Style = $defaultStyle

$ws.Range("D44").Value = "'0.279"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "'  -2.56%  "
$ws.Range("E44").Style = $defaultStyle

$ws.Range("D45").Value = "'2.792.09"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "'  +0.03%  "
$ws.Range("E45").Style = $defaultStyle

$ws.Range("D46").Value = "'0.0354"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "'  -1.01%  "
$ws.Range("E46").Style = $defaultStyle

$ws.Range("D47").Value = "'37.82"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "'  -8.40%  "
$ws.Range("E47").Style = $defaultStyle

$ws.Range("D48").Value = "'129.35"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "'  -1.04%  "
$ws.Range("E48").Style = $defaultStyle

$ws.Range("D50").Value = "'24.17"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "'  -3.26%  "
$ws.Range("E50").Style = $defaultStyle

$ws.Range("E51").Value = "'  -0.73%  "
$ws.Range("E51").Style = $defaultStyle
